# Auto-generated edit script applying numeric corrections to the
# Valefor_Profits sheets (currentAveragePrice / Leve profit columns).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 750.5
$ws.Range("I20").Value = 750.5
$ws.Range("K20").Value = 750.5
$ws.Range("M20").Value = -520.5
$ws.Range("H28").Value = 990.5417
$ws.Range("I28").Value = 1045.5
$ws.Range("J28").Value = 913.6
$ws.Range("K28").Value = 1045.5
$ws.Range("L28").Value = 913.6
$ws.Range("M28").Value = -560.5
$ws.Range("N28").Value = -1883.6
$ws.Range("H35").Value = 750.5
$ws.Range("I35").Value = 750.5
$ws.Range("K35").Value = 750.5
$ws.Range("M35").Value = -371.5
$ws.Range("H62").Value = 1255.0952
$ws.Range("I62").Value = 1171.6471
$ws.Range("J62").Value = 1609.75
$ws.Range("K62").Value = 1171.6471
$ws.Range("L62").Value = 1609.75
$ws.Range("M62").Value = -547.6470999999999
$ws.Range("N62").Value = -2857.75
$ws.Range("H65").Value = 1255.0952
$ws.Range("I65").Value = 1171.6471
$ws.Range("J65").Value = 1609.75
$ws.Range("K65").Value = 5858.2355
$ws.Range("L65").Value = 8048.75
$ws.Range("M65").Value = -2738.2355
$ws.Range("N65").Value = -14288.75
$ws.Range("H98").Value = 42898.668
$ws.Range("I98").Value = 61481.395
$ws.Range("J98").Value = 2016.6666
$ws.Range("K98").Value = 61481.395
$ws.Range("L98").Value = 2016.6666
$ws.Range("M98").Value = -59983.395
$ws.Range("N98").Value = -5012.6666
$ws.Range("H107").Value = 497.26666
$ws.Range("I107").Value = 497.26666
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 497.26666
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 1422.73334
$ws.Range("H111").Value = 1350.25
$ws.Range("I111").Value = 950
$ws.Range("J111").Value = 1483.6666
$ws.Range("K111").Value = 2850
$ws.Range("L111").Value = 4450.9998
$ws.Range("M111").Value = 217
$ws.Range("N111").Value = -10584.9998
$ws.Range("H113").Value = 2176.25
$ws.Range("I113").Value = 1556.4286
$ws.Range("J113").Value = 3044
$ws.Range("K113").Value = 1556.4286
$ws.Range("L113").Value = 3044
$ws.Range("M113").Value = 1697.5714
$ws.Range("N113").Value = -9552
$ws.Range("H116").Value = 1254426.1
$ws.Range("I116").Value = 3335613
$ws.Range("J116").Value = 5714
$ws.Range("K116").Value = 3335613
$ws.Range("L116").Value = 5714
$ws.Range("M116").Value = -3332171
$ws.Range("N116").Value = -12598
$ws.Range("H122").Value = 42898.668
$ws.Range("I122").Value = 61481.395
$ws.Range("J122").Value = 2016.6666
$ws.Range("K122").Value = 184444.185
$ws.Range("L122").Value = 6049.9998
$ws.Range("M122").Value = -181994.185
$ws.Range("N122").Value = -10949.9998
$ws.Range("H125").Value = 4748.6
$ws.Range("J125").Value = 4748.6
$ws.Range("L125").Value = 42737.4
$ws.Range("N125").Value = -47657.4

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5548.5
$ws.Range("I32").Value = 3546.726
$ws.Range("J32").Value = 17559.143
$ws.Range("K32").Value = 3546.726
$ws.Range("L32").Value = 17559.143
$ws.Range("M32").Value = -3259.726
$ws.Range("N32").Value = -18133.143
$ws.Range("H39").Value = 3598.8
$ws.Range("I39").Value = 3598.8
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 3598.8
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -3078.8
$ws.Range("H74").Value = 15152730
$ws.Range("I74").Value = 16667736
$ws.Range("J74").Value = 2666.6667
$ws.Range("K74").Value = 16667736
$ws.Range("L74").Value = 2666.6667
$ws.Range("M74").Value = -16666862
$ws.Range("N74").Value = -4414.6667
$ws.Range("H77").Value = 15152730
$ws.Range("I77").Value = 16667736
$ws.Range("J77").Value = 2666.6667
$ws.Range("K77").Value = 83338680
$ws.Range("L77").Value = 13333.3335
$ws.Range("M77").Value = -83334312
$ws.Range("N77").Value = -22069.3335

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26319032
$ws.Range("I31").Value = 62501244
$ws.Range("J31").Value = 4694.727
$ws.Range("K31").Value = 62501244
$ws.Range("L31").Value = 4694.727
$ws.Range("M31").Value = -62500949
$ws.Range("N31").Value = -5284.727
$ws.Range("H34").Value = 26319032
$ws.Range("I34").Value = 62501244
$ws.Range("J34").Value = 4694.727
$ws.Range("K34").Value = 62501244
$ws.Range("L34").Value = 4694.727
$ws.Range("M34").Value = -62501042
$ws.Range("N34").Value = -5098.727
$ws.Range("H63").Value = 30271
$ws.Range("J63").Value = 30271
$ws.Range("L63").Value = 30271
$ws.Range("N63").Value = -31643
$ws.Range("H66").Value = 30271
$ws.Range("J66").Value = 30271
$ws.Range("L66").Value = 90813
$ws.Range("N66").Value = -97677
$ws.Range("H122").Value = 4520.4062
$ws.Range("I122").Value = 5463.6523
$ws.Range("J122").Value = 2109.889
$ws.Range("K122").Value = 16390.9569
$ws.Range("L122").Value = 6329.667
$ws.Range("M122").Value = -13940.9569
$ws.Range("N122").Value = -11229.667
$ws.Range("H141").Value = 32339.8
$ws.Range("J141").Value = 32339.8
$ws.Range("L141").Value = 32339.8
$ws.Range("N141").Value = -42699.8

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 4800
$ws.Range("J48").Value = 4800
$ws.Range("L48").Value = 14400
$ws.Range("N48").Value = -14900
$ws.Range("H58").Value = 3375
$ws.Range("J58").Value = 4550
$ws.Range("L58").Value = 13650
$ws.Range("N58").Value = -13906

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 27500
$ws.Range("I48").Value = 5000
$ws.Range("K48").Value = 5000
$ws.Range("M48").Value = -4515
$ws.Range("H80").Value = 168551
$ws.Range("I80").Value = 2233.3333
$ws.Range("J80").Value = 334868.66
$ws.Range("K80").Value = 2233.3333
$ws.Range("L80").Value = 334868.66
$ws.Range("M80").Value = -1235.3333
$ws.Range("N80").Value = -336864.66
$ws.Range("H83").Value = 168551
$ws.Range("I83").Value = 2233.3333
$ws.Range("J83").Value = 334868.66
$ws.Range("K83").Value = 11166.6665
$ws.Range("L83").Value = 1674343.3
$ws.Range("M83").Value = -6174.666499999999
$ws.Range("N83").Value = -1684327.3

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2273.6667
$ws.Range("I61").Value = 1346.3846
$ws.Range("J61").Value = 3369.5454
$ws.Range("K61").Value = 1346.3846
$ws.Range("L61").Value = 3369.5454
$ws.Range("M61").Value = -1144.3846
$ws.Range("N61").Value = -3773.5454
$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -9376
$ws.Range("H64").Value = 29214.334
$ws.Range("I64").Value = 41568
$ws.Range("J64").Value = 23037.5
$ws.Range("K64").Value = 41568
$ws.Range("L64").Value = 23037.5
$ws.Range("M64").Value = -41343
$ws.Range("N64").Value = -23487.5
$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -26880
$ws.Range("H67").Value = 29214.334
$ws.Range("I67").Value = 41568
$ws.Range("J67").Value = 23037.5
$ws.Range("K67").Value = 41568
$ws.Range("L67").Value = 23037.5
$ws.Range("M67").Value = -40788
$ws.Range("N67").Value = -24597.5
$ws.Range("H113").Value = 2273.6667
$ws.Range("I113").Value = 1346.3846
$ws.Range("J113").Value = 3369.5454
$ws.Range("K113").Value = 1346.3846
$ws.Range("L113").Value = 3369.5454
$ws.Range("M113").Value = 823.6153999999999
$ws.Range("N113").Value = -7709.5454

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 31299.334
$ws.Range("J54").Value = 31299.334
$ws.Range("L54").Value = 31299.334
$ws.Range("N54").Value = -32339.334
$ws.Range("H81").Value = 36303.438
$ws.Range("I81").Value = 15421.571
$ws.Range("J81").Value = 42150.36
$ws.Range("K81").Value = 30843.142
$ws.Range("L81").Value = 84300.72
$ws.Range("M81").Value = -29782.142
$ws.Range("N81").Value = -86422.72
$ws.Range("H84").Value = 36303.438
$ws.Range("I84").Value = 15421.571
$ws.Range("J84").Value = 42150.36
$ws.Range("K84").Value = 154215.71
$ws.Range("L84").Value = 421503.6
$ws.Range("M84").Value = -148911.71
$ws.Range("N84").Value = -432111.6
$ws.Range("H107").Value = 3944.7778
$ws.Range("I107").Value = 223.17647
$ws.Range("J107").Value = 10271.5
$ws.Range("K107").Value = 669.52941
$ws.Range("L107").Value = 30814.5
$ws.Range("M107").Value = 1250.47059
$ws.Range("N107").Value = -34654.5
$ws.Range("H113").Value = 677.3333
$ws.Range("I113").Value = 270.4
$ws.Range("J113").Value = 968
$ws.Range("K113").Value = 811.1999999999999
$ws.Range("L113").Value = 2904
$ws.Range("M113").Value = 1358.8
$ws.Range("N113").Value = -7244
$ws.Range("H122").Value = 1343.9286
$ws.Range("I122").Value = 1081.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3244.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -794.5
$ws.Range("N122").Value = -10900
$ws.Range("H126").Value = 2194.3572
$ws.Range("I126").Value = 1320
$ws.Range("J126").Value = 2680.111
$ws.Range("K126").Value = 3960
$ws.Range("L126").Value = 8040.333
$ws.Range("M126").Value = -1490
$ws.Range("N126").Value = -12980.333

Write-Host "Applied all Valefor_Profits numeric updates."